# Apply cryptos list update (values + percentages), matching the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.475.24"
$ws.Cells.Item(2, 5).Value = "  +0.98%  "

$ws.Cells.Item(3, 4).Value = "1.920.57"
$ws.Cells.Item(3, 5).Value = "  +1.80%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.009"
$ws.Cells.Item(4, 5).Value = "  +0.85%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "325.65"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.007"
$ws.Cells.Item(6, 5).Value = "  +0.66%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4826"
$ws.Cells.Item(7, 5).Value = "  +2.60%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4087"
$ws.Cells.Item(8, 5).Value = "  +1.64%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.08168"
$ws.Cells.Item(9, 5).Value = "  +2.01%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.023"
$ws.Cells.Item(10, 5).Value = "  +2.91%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "23.52"
$ws.Cells.Item(11, 5).Value = "  +5.06%  "

$ws.Cells.Item(12, 4).Value = "1.907.47"
$ws.Cells.Item(12, 5).Value = "  +1.55%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.045"
$ws.Cells.Item(13, 5).Value = "  +2.79%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.235"
$ws.Cells.Item(14, 5).Value = "  +2.84%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "91.33"
$ws.Cells.Item(15, 5).Value = "  +2.87%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.009"
$ws.Cells.Item(16, 5).Value = "  +0.82%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.06784"
$ws.Cells.Item(17, 5).Value = "  +2.50%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.00001040"
$ws.Cells.Item(18, 5).Value = "  +1.61%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "17.77"
$ws.Cells.Item(19, 5).Value = "  +1.93%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "1.007"
$ws.Cells.Item(20, 5).Value = "  +0.60%  "

$ws.Cells.Item(21, 4).Value = "29.511.22"
$ws.Cells.Item(21, 5).Value = "  +1.13%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.630"
$ws.Cells.Item(22, 5).Value = "  +2.63%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "11.76"
$ws.Cells.Item(23, 5).Value = "  +0.53%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.188"
$ws.Cells.Item(24, 5).Value = "  +0.53%  "

$ws.Cells.Item(25, 4).Value = "2.161.36"
$ws.Cells.Item(25, 5).Value = "  +2.11%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "6.731"
$ws.Cells.Item(26, 5).Value = "  +11.78%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "156.41"
$ws.Cells.Item(27, 5).Value = "  +1.00%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "20.08"
$ws.Cells.Item(28, 5).Value = "  +2.40%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.116"
$ws.Cells.Item(29, 5).Value = "  +2.06%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "120.42"
$ws.Cells.Item(30, 5).Value = "  +2.68%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.029"
$ws.Cells.Item(31, 5).Value = "  -0.24%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.09564"
$ws.Cells.Item(32, 5).Value = "  +1.49%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.521"
$ws.Cells.Item(33, 5).Value = "  +3.22%  "

$ws.Cells.Item(34, 5).Value = "  +0.89%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.389"
$ws.Cells.Item(35, 5).Value = "  +0.73%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.02280"
$ws.Cells.Item(36, 5).Value = "  +2.49%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.06134"
$ws.Cells.Item(37, 5).Value = "  +1.13%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.180"
$ws.Cells.Item(38, 5).Value = "  +0.60%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.5973"
$ws.Cells.Item(39, 5).Value = "  +2.82%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "8.034"
$ws.Cells.Item(40, 5).Value = "  +0.23%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "10.77"
$ws.Cells.Item(41, 5).Value = "  +7.69%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.007"
$ws.Cells.Item(42, 5).Value = "  +0.79%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.1860"
$ws.Cells.Item(43, 5).Value = "  +1.89%  "

$ws.Cells.Item(44, 2).Value = "RenderToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.417"
$ws.Cells.Item(44, 5).Value = "  -2.64%  "

$ws.Cells.Item(45, 2).Value = "WEMIXToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.249"
$ws.Cells.Item(45, 5).Value = "  -1.90%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.07609"
$ws.Cells.Item(46, 5).Value = "  -1.02%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "12.45"
$ws.Cells.Item(47, 5).Value = "  +3.06%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.5585"
$ws.Cells.Item(48, 5).Value = "  +2.11%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.956"
$ws.Cells.Item(49, 5).Value = "  +2.99%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "117.20"
$ws.Cells.Item(50, 5).Value = "  +3.26%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.440"
$ws.Cells.Item(51, 5).Value = "  +5.03%  "
